$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 2877
$wsExhibition.Range("F12").Value = 42
$wsExhibition.Range("F14").Value = 947

# Sheet "全部类型" (All Types) - same underlying rows, update corresponding values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 2877
$wsAll.Range("F14").Value = 42
$wsAll.Range("F16").Value = 947
